$p = $ppt.ActivePresentation

# Delete the old last slide (22, "Kocka na uklidneni") - its content/position
# is being folded into slide 21 instead, and slide 21 survives as the deck's
# final slide.
$p.Slides.Item(22).Delete()

# Slide 21 ("Konec") keeps its title, but the illustration is enlarged and
# the caption textbox is repurposed with new text, moved near the bottom.
$s21 = $p.Slides.Item(21)

$pic = $s21.Shapes.Item(2)
$pic.Left = 2483027 / 12700
$pic.Top = 1731963 / 12700
$pic.Width = 7216421 / 12700
$pic.Height = 4059237 / 12700

$caption = $s21.Shapes.Item(3)
$caption.Left = 913795 / 12700
$caption.Top = 5692923 / 12700
$caption.TextFrame.TextRange.Text = "Kotě na uklidnění"
$caption.TextFrame.TextRange.Font.Size = 28
